$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("B2").Value = -0.4999999999999858
$ws.Range("C2").Value = -0.6999999999999886
$ws.Range("D2").Value = 1.000000000000014
$ws.Range("E2").Value = 0.5000000000000284
$ws.Range("F2").Value = -2.799999999999969
$ws.Range("G2").Value = 0.00000000000002842170943040401
$ws.Range("H2").Value = -0.1999999999999744
$ws.Range("B3").Value = -0.2000000000000028
$ws.Range("C3").Value = 1.5
$ws.Range("D3").Value = 1.000000000000014
$ws.Range("E3").Value = -2.299999999999983
$ws.Range("F3").Value = 0.5000000000000142
$ws.Range("G3").Value = 0.3000000000000114
$ws.Range("B4").Value = 1.700000000000003
$ws.Range("C4").Value = 1.200000000000017
$ws.Range("D4").Value = -2.09999999999998
$ws.Range("E4").Value = 0.7000000000000171
$ws.Range("F4").Value = 0.5000000000000142
$ws.Range("G4").Value = 0.8000000000000114
$ws.Range("H4").Value = 0.5000000000000142
$ws.Range("I4").Value = 0.6000000000000227
$ws.Range("J4").Value = 0.6000000000000227
$ws.Range("B5").Value = -0.4999999999999858
$ws.Range("C5").Value = -3.799999999999983
$ws.Range("D5").Value = -0.9999999999999858
$ws.Range("E5").Value = -1.199999999999989
$ws.Range("F5").Value = -0.8999999999999915
$ws.Range("G5").Value = -1.199999999999989
$ws.Range("H5").Value = -1.09999999999998
$ws.Range("I5").Value = -1.09999999999998
$ws.Range("B6").Value = -3.299999999999997
$ws.Range("C6").Value = -0.5
$ws.Range("D6").Value = -0.7000000000000028
$ws.Range("E6").Value = -0.4000000000000057
$ws.Range("F6").Value = -0.7000000000000028
$ws.Range("G6").Value = -0.5999999999999943
$ws.Range("H6").Value = -0.5999999999999943
$ws.Range("B7").Value = 2.799999999999997
$ws.Range("C7").Value = 2.599999999999994
$ws.Range("D7").Value = 2.899999999999991
$ws.Range("E7").Value = 2.599999999999994
$ws.Range("F7").Value = 2.700000000000003
$ws.Range("G7").Value = 2.700000000000003
$ws.Range("B8").Value = -0.2000000000000028
$ws.Range("C8").Value = 0.09999999999999432
$ws.Range("D8").Value = -0.2000000000000028
$ws.Range("E8").Value = -0.09999999999999432
$ws.Range("F8").Value = -0.09999999999999432
$ws.Range("G8").Value = 0
$ws.Range("H8").Value = -0.4999999999999858
$ws.Range("I8").Value = -0.4000000000000199
$ws.Range("B9").Value = 0.2999999999999971
$ws.Range("C9").Value = 0
$ws.Range("D9").Value = 0.1000000000000085
$ws.Range("E9").Value = 0.1000000000000085
$ws.Range("F9").Value = 0.2000000000000028
$ws.Range("G9").Value = -0.299999999999983
$ws.Range("H9").Value = -0.2000000000000171
$ws.Range("B10").Value = -0.2999999999999971
$ws.Range("C10").Value = -0.1999999999999886
$ws.Range("D10").Value = -0.1999999999999886
$ws.Range("E10").Value = -0.09999999999999432
$ws.Range("F10").Value = -0.5999999999999801
$ws.Range("G10").Value = -0.5000000000000142
$ws.Range("B11").Value = 0.1000000000000085
$ws.Range("C11").Value = 0.1000000000000085
$ws.Range("D11").Value = 0.2000000000000028
$ws.Range("E11").Value = -0.299999999999983
$ws.Range("F11").Value = -0.2000000000000171
$ws.Range("B12").Value = 0
$ws.Range("C12").Value = 0.09999999999999432
$ws.Range("D12").Value = -0.3999999999999915
$ws.Range("E12").Value = -0.3000000000000256
$ws.Range("B13").Value = 0.09999999999999432
$ws.Range("C13").Value = -0.3999999999999915
$ws.Range("D13").Value = -0.3000000000000256
$ws.Range("B14").Value = -0.4999999999999858
$ws.Range("C14").Value = -0.4000000000000199
$ws.Range("B15").Value = 0.09999999999996589
